$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# ---------------------------------------------------------------------
# 1) Clone cell formatting (styles only, no values) from the fully
#    populated template row 113 onto every cell we are about to touch,
#    so new cells land on the same style indices (s="...") as the rest
#    of the table, exactly like the target file shows.
# ---------------------------------------------------------------------

# F114:F118 - column F doesn't exist yet on those rows (style 1)
$ws.Range("F113").Copy()
$ws.Range("F114:F118").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 118 needs the full A:L template (A-C/G get real values below,
# D/H/I/J/K/L stay "blank but styled").
$ws.Range("A113:L113").Copy()
$ws.Range("A118:L118").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 119-124: clone the A:D and G:L formatting template from row 113
# (E is handled separately below via the new shared formula) - this
# creates the "blank" s=6/7/8/8 ... s=18/13/13/13/13/13 cells with no
# values, matching the target exactly (no F column on these rows).
$ws.Range("A113:D113").Copy()
$ws.Range("A119:D124").PasteSpecial(-4122)
$ws.Range("G113:L113").Copy()
$ws.Range("G119:L124").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Populate new shared-string text in the exact order the target
#    workbook first introduces each one, so the appended sharedStrings
#    entries land on the expected indices (141-145).
# ---------------------------------------------------------------------
$ws.Range("G114").Value = "Gestion des articles"         # -> new string #141
$ws.Range("G117").Value = "Pop-up de confirmation"        # -> new string #142
$ws.Range("H115").Value = "Suppression"                   # -> new string #143
$ws.Range("K117").Value = "https://www.tutorialrepublic.com/faq/how-to-get-the-current-url-with-javascript.php`nhttps://www.w3schools.com/howto/howto_js_redirect_webpage.asp`nhttps://www.w3schools.com/tags/att_script_src.asp`nhttps://www.w3schools.com/js/js_window_location.asp`nhttps://www.w3schools.com/howto/tryit.asp?filename=tryhow_js_redirect_webpage`nhttps://developer.mozilla.org/fr/docs/Web/API/window/location`n"   # -> new string #144
$ws.Range("K115").Value = "https://www.w3schools.com/sql/sql_delete.asp"   # -> new string #145
$ws.Range("K115").Style = "Normal"   # target drops K115's style back to default (no s="13")

# ---------------------------------------------------------------------
# 3) Fill in the rest of row 114 - Gestion des articles / Modification
# ---------------------------------------------------------------------
$ws.Range("A114").Value = 44343
$ws.Range("B114").Value = 4
$ws.Range("C114").Value = 0.5625
$ws.Range("D114").Value = 0.56944444444444442
$ws.Range("F114").Value = "Réalisation"
$ws.Range("H114").Value = "Modification"

# ---------------------------------------------------------------------
# 4) Row 115 - Gestion des œuvres / Suppression
# ---------------------------------------------------------------------
$ws.Range("A115").Value = 44343
$ws.Range("B115").Value = 4
$ws.Range("C115").Value = 0.56944444444444442
$ws.Range("D115").Value = 0.58333333333333337
$ws.Range("F115").Value = "Réalisation"
$ws.Range("G115").Value = "Gestion des œuvres"

# ---------------------------------------------------------------------
# 5) Row 116 - Gestion des articles / Suppression
# ---------------------------------------------------------------------
$ws.Range("A116").Value = 44343
$ws.Range("B116").Value = 4
$ws.Range("C116").Value = 0.58333333333333337
$ws.Range("D116").Value = 0.60416666666666663
$ws.Range("F116").Value = "Réalisation"
$ws.Range("G116").Value = "Gestion des articles"
$ws.Range("H116").Value = "Suppression"

# ---------------------------------------------------------------------
# 6) Row 117 - Pop-up de confirmation (tall wrapped row)
# ---------------------------------------------------------------------
$ws.Range("A117").Value = 44343
$ws.Range("B117").Value = 4
$ws.Range("C117").Value = 0.60416666666666663
$ws.Range("D117").Value = 0.62847222222222221
$ws.Range("F117").Value = "Réalisation"
$ws.Range("A117").EntireRow.RowHeight = 195

# ---------------------------------------------------------------------
# 7) Row 118 - Pop-up de confirmation, new shared-formula block starts
#    here (E118:E124), D118 stays empty so E118 = 0 - C118.
# ---------------------------------------------------------------------
$ws.Range("A118").Value = 44343
$ws.Range("B118").Value = 4
$ws.Range("C118").Value = 0.63888888888888895
$ws.Range("F118").Value = "Réalisation"
$ws.Range("G118").Value = "Pop-up de confirmation"

# New shared formula group si="9" covering E118:E124
$ws.Range("E118:E124").Formula = "=D118-C118"

Write-Host "Edits applied"
